$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so they stay text (matching the original inlineStr type) instead of
# being auto-converted to numbers by Excel.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range('D2').Value = '30.990.75'
$ws.Range('E2').Value = '  +2.89%  '
$ws.Range('D3').Value = '2.117.57'
$ws.Range('E3').Value = '  +10.32%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '334.37'
$ws.Range('E5').Value = '  +4.70%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').Value = '0.5337'
$ws.Range('E7').Value = '  +5.19%  '
$ws.Range('D8').Value = '0.4404'
$ws.Range('E8').Value = '  +8.38%  '
$ws.Range('D9').Value = '0.09068'
$ws.Range('E9').Value = '  +8.78%  '
$ws.Range('D10').Value = '46.35'
$ws.Range('E10').Value = '  +10.13%  '
$ws.Range('E12').Value = '  +4.43%  '
$ws.Range('D13').Value = '2.119.46'
$ws.Range('E13').Value = '  +10.42%  '
$ws.Range('D14').Value = '6.802'
$ws.Range('E14').Value = '  +5.99%  '
$ws.Range('D15').Value = '7.830'
$ws.Range('E15').Value = '  +7.92%  '
$ws.Range('D16').Value = '97.87'
$ws.Range('E16').Value = '  +5.68%  '
$ws.Range('D17').Value = '1.002'
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = '0.00001140'
$ws.Range('E18').Value = '  +4.04%  '
$ws.Range('D19').Value = '0.06683'
$ws.Range('E19').Value = '  +2.85%  '
$ws.Range('E20').Value = '  +4.01%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '6.399'
$ws.Range('E22').Value = '  +7.60%  '
$ws.Range('D23').Value = '31.067.76'
$ws.Range('E23').Value = '  +3.14%  '
$ws.Range('D24').Value = '12.24'
$ws.Range('E24').Value = '  +7.82%  '
$ws.Range('D25').Value = '2.371.87'
$ws.Range('E26').Value = '  +3.60%  '
$ws.Range('E27').Value = '  +4.98%  '
$ws.Range('D28').Value = '2.564'
$ws.Range('E28').Value = '  +13.19%  '
$ws.Range('D29').Value = '163.70'
$ws.Range('E29').Value = '  +0.76%  '
$ws.Range('D30').Value = '133.94'
$ws.Range('E30').Value = '  +3.84%  '
$ws.Range('D31').Value = '1.181'
$ws.Range('E31').Value = '  +3.96%  '
$ws.Range('D32').Value = '0.1078'
$ws.Range('E32').Value = '  +2.96%  '
$ws.Range('D33').Value = '6.269'
$ws.Range('E33').Value = '  +5.45%  '
$ws.Range('D34').Value = '4.009'
$ws.Range('E34').Value = '  +5.63%  '
$ws.Range('D35').Value = '1.534'
$ws.Range('E35').Value = '  +25.39%  '
$ws.Range('E36').Value = '  +7.43%  '
$ws.Range('E37').Value = '  +14.91%  '
$ws.Range('D38').Value = '5.574'
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.06767'
$ws.Range('E39').Value = '  +5.27%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '9.557'
$ws.Range('E40').Value = '  +11.18%  '
$ws.Range('D41').Value = '0.2287'
$ws.Range('E41').Value = '  +6.57%  '
$ws.Range('D42').Value = '0.6892'
$ws.Range('E42').Value = '  +6.60%  '
$ws.Range('D43').Value = '1.261'
$ws.Range('E43').Value = '  +3.91%  '
$ws.Range('D44').Value = '0.6516'
$ws.Range('E44').Value = '  +7.75%  '
$ws.Range('D45').Value = '14.12'
$ws.Range('E45').Value = '  +5.60%  '
$ws.Range('D46').Value = '1.001'
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('D47').Value = '2.265'
$ws.Range('E47').Value = '  +4.15%  '
$ws.Range('D48').Value = '3.679'
$ws.Range('E48').Value = '  +1.52%  '
$ws.Range('E49').Value = '  +6.03%  '
$ws.Range('D50').Value = '83.32'
$ws.Range('E50').Value = '  +6.73%  '
$ws.Range('D51').Value = '120.07'
$ws.Range('E51').Value = '  -1.81%  '
